$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1112.5
$ws.Range("J2").Value = 725
$ws.Range("L2").Value = 725
$ws.Range("N2").Value = -951
$ws.Range("H15").Value = 537.375
$ws.Range("I15").Value = 537.375
$ws.Range("K15").Value = 1612.125
$ws.Range("M15").Value = -1443.125
$ws.Range("H70").Value = 3499.9
$ws.Range("I70").Value = 3388.889
$ws.Range("K70").Value = 10166.667
$ws.Range("M70").Value = -9896.667000000001
$ws.Range("H73").Value = 3499.9
$ws.Range("I73").Value = 3388.889
$ws.Range("K73").Value = 10166.667
$ws.Range("M73").Value = -9230.667000000001
$ws.Range("H88").Value = 2053
$ws.Range("J88").Value = 2159.3
$ws.Range("L88").Value = 2159.3
$ws.Range("N88").Value = -2971.3
$ws.Range("H91").Value = 2053
$ws.Range("J91").Value = 2159.3
$ws.Range("L91").Value = 2159.3
$ws.Range("N91").Value = -4967.3
$ws.Range("H116").Value = 5995
$ws.Range("I116").Value = 4990
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 4990
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -1548
$ws.Range("N116").Value = -13884
$ws.Range("H138").Value = 4251.081
$ws.Range("J138").Value = 4349.759
$ws.Range("L138").Value = 13049.277
$ws.Range("N138").Value = -23329.277
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1389.3334
$ws.Range("I45").Value = 1251.5
$ws.Range("J45").Value = 1665
$ws.Range("K45").Value = 1251.5
$ws.Range("L45").Value = 1665
$ws.Range("M45").Value = -874.5
$ws.Range("N45").Value = -2419
$ws.Range("H96").Value = 18348.8
$ws.Range("J96").Value = 18348.8
$ws.Range("L96").Value = 18348.8
$ws.Range("N96").Value = -23840.8
$ws.Range("H119").Value = 33125
$ws.Range("J119").Value = 33125
$ws.Range("L119").Value = 33125
$ws.Range("N119").Value = -42801
$ws.Range("H122").Value = 4866.0625
$ws.Range("I122").Value = 2532.4546
$ws.Range("K122").Value = 7597.3638
$ws.Range("M122").Value = -5147.3638

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 221581.6
$ws.Range("I94").Value = 221581.6
$ws.Range("K94").Value = 221581.6
$ws.Range("M94").Value = -221130.6
$ws.Range("H99").Value = 2599.8572
$ws.Range("I99").Value = 2900
$ws.Range("J99").Value = 799
$ws.Range("K99").Value = 2900
$ws.Range("L99").Value = 799
$ws.Range("M99").Value = -1402
$ws.Range("N99").Value = -3795

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8663.223
$ws.Range("J31").Value = 9162.666999999999
$ws.Range("L31").Value = 9162.666999999999
$ws.Range("N31").Value = -9752.666999999999
$ws.Range("H34").Value = 8663.223
$ws.Range("J34").Value = 9162.666999999999
$ws.Range("L34").Value = 9162.666999999999
$ws.Range("N34").Value = -9566.666999999999
$ws.Range("H99").Value = 1115034.5
$ws.Range("I99").Value = 629413.75
$ws.Range("K99").Value = 629413.75
$ws.Range("M99").Value = -627915.75
$ws.Range("H105").Value = 2649.6428
$ws.Range("I105").Value = 2626.9092
$ws.Range("K105").Value = 2626.9092
$ws.Range("M105").Value = -879.9092000000001
$ws.Range("H126").Value = 1115034.5
$ws.Range("I126").Value = 629413.75
$ws.Range("K126").Value = 1888241.25
$ws.Range("M126").Value = -1885771.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3041.5
$ws.Range("J58").Value = 3500
$ws.Range("L58").Value = 10500
$ws.Range("N58").Value = -10756
$ws.Range("H86").Value = 183
$ws.Range("J86").Value = 199.5
$ws.Range("L86").Value = 598.5
$ws.Range("N86").Value = -2970.5
$ws.Range("H89").Value = 183
$ws.Range("J89").Value = 199.5
$ws.Range("L89").Value = 1795.5
$ws.Range("N89").Value = -13651.5
$ws.Range("H113").Value = 685.6
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 685.6
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2056.8
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6396.8

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 12500
$ws.Range("J53").Value = 15000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -16262
$ws.Range("H80").Value = 28250
$ws.Range("I80").Value = 28250
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 28250
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -27252
$ws.Range("H83").Value = 28250
$ws.Range("I83").Value = 28250
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 141250
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -136258
$ws.Range("H122").Value = 2219.2
$ws.Range("I122").Value = 1498.6666
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 4495.9998
$ws.Range("M122").Value = -2045.9998
$ws.Range("N122").Value = -14800

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2478.8
$ws.Range("I7").Value = 2478.8
$ws.Range("K7").Value = 2478.8
$ws.Range("M7").Value = -2366.8
$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 20000
$ws.Range("M57").Value = -19434
$ws.Range("H74").Value = 65739
$ws.Range("J74").Value = 65739
$ws.Range("L74").Value = 65739
$ws.Range("N74").Value = -67735
$ws.Range("H77").Value = 65739
$ws.Range("J77").Value = 65739
$ws.Range("L77").Value = 197217
$ws.Range("N77").Value = -207201
$ws.Range("H126").Value = 2478.8
$ws.Range("I126").Value = 2478.8
$ws.Range("K126").Value = 7436.400000000001
$ws.Range("M126").Value = -4966.400000000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H132").Value = 921.2
$ws.Range("I132").Value = 921.2
$ws.Range("K132").Value = 2763.6
$ws.Range("M132").Value = -233.6000000000004
